$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 10 - "Content Placeholder 2" (Shapes.Item(2))
# ---------------------------------------------------------------------------
$s10 = $p.Slides.Item(10)
$tr10 = $s10.Shapes.Item(2).TextFrame.TextRange

# Para 2: "\tNeutral/Disagree" -> "Disagree"
$tr10.Paragraphs(2).Text = "Disagree"

# Para 5: "\tAgree" -> " " + "      Agree" (tab becomes a single space, label gets
# leading spaces) -- set via two runs using InsertBefore/replace.
$tr10.Paragraphs(5).Text = "       Agree"
$tr10.Paragraphs(5).Characters(1, 1).Text = " "

# Remove paragraphs 7-9 ("The organization of information on the screen was
# clear", "	Neutral/Disagree", and the blank line that followed them).
$para7 = $tr10.Paragraphs(7)
$para9 = $tr10.Paragraphs(9)
$startPos = $para7.Start
$len = ($para9.Start + $para9.Length) - $startPos
$tr10.Characters($startPos, $len).Delete()

# "Agree/Neutral" -> "Agree" (now paragraph 8 after the deletion above)
$tr10.Paragraphs(8).Text = "Agree"

# Remove the now-orphaned empty numbered paragraph that used to sit right
# after "Agree/Neutral".
$paraBlank = $tr10.Paragraphs(9)
$tr10.Characters($paraBlank.Start, $paraBlank.Length).Delete()

# ---------------------------------------------------------------------------
# Slide 9 - "Content Placeholder 2" (Shapes.Item(2))
# ---------------------------------------------------------------------------
$s9 = $p.Slides.Item(9)
$tr9 = $s9.Shapes.Item(2).TextFrame.TextRange

# Para 2: "\tEasy/Moderate" -> "Neutral/Agree"
$tr9.Paragraphs(2).Text = "Neutral/Agree"

# Para 5: "\tAgree" -> "Agree"
$tr9.Paragraphs(5).Text = "Agree"

# Para 8: "I prefer using Smart-Waiter over traditional sense" -> split into
# two runs ("I prefer using Smart-Waiter over traditional " + "sense").
$tr9.Paragraphs(8).Text = "I prefer using Smart-Waiter over traditional sense"

# Para 9: "\tNo Preference/Agree" -> "Agree"
$tr9.Paragraphs(9).Text = "Agree"

# Para 12: "The interface of the system was pleasant" -> split into two runs
# ("The interface of the system was " + "pleasant").
$tr9.Paragraphs(12).Text = "The interface of the system was pleasant"

# Para 13: "\tDisagree" -> "Disagree"
$tr9.Paragraphs(13).Text = "Disagree"
